# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" / "演出" / "本地生活" / "全部类型" sheets, per the latest
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F6").Value  = 922
$wsExhibit.Range("F7").Value  = 468
$wsExhibit.Range("F9").Value  = 2209
$wsExhibit.Range("F12").Value = 123
$wsExhibit.Range("F13").Value = 1106
$wsExhibit.Range("F14").Value = 186
$wsExhibit.Range("F15").Value = 2226
$wsExhibit.Range("F16").Value = 684
$wsExhibit.Range("F17").Value = 13585
$wsExhibit.Range("F19").Value = 1289
$wsExhibit.Range("F20").Value = 47
$wsExhibit.Range("F25").Value = 64
$wsExhibit.Range("F27").Value = 275
$wsExhibit.Range("F30").Value = 29

# Sheet 2: 演出
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F5").Value  = 20
$wsShow.Range("F17").Value = 21
$wsShow.Range("F19").Value = 2

# Sheet 3: 本地生活
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Range("F2").Value = 5721

# Sheet 4: 全部类型
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F8").Value  = 922
$wsAll.Range("F10").Value = 468
$wsAll.Range("F11").Value = 20
$wsAll.Range("F12").Value = 2209
$wsAll.Range("F16").Value = 123
$wsAll.Range("F18").Value = 1106
$wsAll.Range("F20").Value = 186
$wsAll.Range("F23").Value = 2226
$wsAll.Range("F24").Value = 684
$wsAll.Range("F27").Value = 1289
$wsAll.Range("F28").Value = 47
$wsAll.Range("F33").Value = 64
$wsAll.Range("F38").Value = 275
$wsAll.Range("F41").Value = 2
$wsAll.Range("F49").Value = 29
